# Updated the input files in sixteen_tests to have strain_log2_expression
# instead of just strain and then ran the files and saved the outputs in
# sixteen_tests_output.
#
# Concretely, for this workbook:
#   - rename worksheet "wt"    -> "wt_log2_expression"
#   - rename worksheet "dcin5" -> "dcin5_log2_expression"
#   - the selection on the (renamed) wt_log2_expression sheet moves from
#     N1:N5 to the single cell G22 (as left behind after someone clicked
#     around on that sheet while renaming it)

$wb = $excel.ActiveWorkbook

# Remember whichever sheet is active right now so we can restore it - the
# workbook's active tab (activeTab) must stay unchanged; only the wt sheet's
# own remembered selection should move to G22.
$originalActive = $wb.ActiveSheet

$wsWt = $wb.Worksheets.Item("wt")
$wsDcin5 = $wb.Worksheets.Item("dcin5")

$wsWt.Name = "wt_log2_expression"
$wsDcin5.Name = "dcin5_log2_expression"

# Touch the wt sheet so its stored selection becomes G22, then flip back to
# whatever sheet was active before so the workbook-level active tab is
# unaffected.
$wsWt.Select() | Out-Null
$wsWt.Range("G22").Select() | Out-Null

$originalActive.Select() | Out-Null
